$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Rename the merge fields in the first heading:
#    {table.monographNumber}: {table.agent}
#      -> {table.volumeNumber}: {table.monographAgent}
#    Done as two independent field renames (monographNumber ->
#    volumeNumber, agent -> monographAgent) to mirror the underlying
#    template-field rename that produced this edit.
# ------------------------------------------------------------------
$d.Content.Find.Execute("{table.monographNumber}", $true, $false, $false, $false, $false, $true, 1, $false, "{table.volumeNumber}", 2)
$d.Content.Find.Execute("{table.agent}", $true, $false, $false, $false, $false, $true, 1, $false, "{table.monographAgent}", 2)

# ------------------------------------------------------------------
# 2) Add two new character styles, ListLabel8 / ListLabel9, mirroring
#    the existing ListLabel6 / ListLabel7 styles already present.
# ------------------------------------------------------------------
$listLabel8 = $d.Styles.Add("ListLabel8", 2)
$listLabel8.NameLocal = "ListLabel 8"
$listLabel8.Font.NameBi = "Symbol"

$listLabel9 = $d.Styles.Add("ListLabel9", 2)
$listLabel9.NameLocal = "ListLabel 9"
$listLabel9.Font.NameBi = "OpenSymbol"

# ------------------------------------------------------------------
# 3) Left-align the Title and Subtitle paragraph styles.
# ------------------------------------------------------------------
$title = $d.Styles("Title")
$title.ParagraphFormat.Alignment = 0

$subtitle = $d.Styles("Subtitle")
$subtitle.ParagraphFormat.Alignment = 0
